$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "pt_max" column (column F, which held the constant value 50 for every
# row) was removed from the sheet. Deleting the entire column shifts every
# column to its right (G:M -> F:L) one position to the left, which also
# updates the dependent formula references (e.g. ABS(0.03*H2) -> ABS(0.03*G2))
# and the shared-strings table (the now-unused "pt_max" string is dropped).
$ws.Columns("F:F").Delete()

# Reflect the new active selection recorded in the saved sheet view.
$ws.Range("D8").Select()
